$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") values change from 45233 to 45243 for rows 2-70
for ($r = 2; $r -le 70; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45233) {
        $cell.Value2 = 45243
    }
}
